# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Column D cells are apostrophe-prefixed so Excel stores the numeric-
# looking price strings (e.g. "42.47", "35.026.81") as text, matching
# the workbook's original inline-string storage for that column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.026.81"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'1.847.74"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'234.66"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'42.47"
$ws.Range("E8").Value = "  +7.38%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'2.117.31"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "'1.857.93"
$ws.Range("E13").Value = "  +2.77%  "
$ws.Range("D14").Value = "'11.36"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").Value = "'34.999.73"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "'69.98"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'240.73"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'4.78"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'2.27"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").Value = "'170.82"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  +21.04%  "
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'0.0555"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").Value = "'3.96"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +23.68%  "
$ws.Range("E35").Value = "  +11.68%  "
$ws.Range("E36").Value = "  +10.00%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +10.55%  "
$ws.Range("D39").Value = "'91.36"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("E40").Value = "  +4.66%  "
$ws.Range("D41").Value = "'1.343.74"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "'14.97"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("B43").Value = "Gas"
$ws.Range("C43").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D43").Value = "'12.97"
$ws.Range("E43").Value = "  +86.69%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.33"
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D47").Value = "'6.38"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("D48").Value = "'0.0533"
$ws.Range("E48").Value = "  +3.80%  "
$ws.Range("D49").Value = "'2.028.06"
$ws.Range("E50").Value = "  +16.16%  "
$ws.Range("D51").Value = "'0.0677"
$ws.Range("E51").Value = "  +1.50%  "
